$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.11739948516318
$ws.Range("C2").Value = 12.40874363839383
$ws.Range("D2").Value = 4.157876949939991
$ws.Range("F2").Value = 19.61730800570504
$ws.Range("G2").Value = 20.72052726873432
$ws.Range("H2").Value = 12.27752810465199
$ws.Range("I2").Value = 18.850489163605
$ws.Range("L2").Value = 10.62282586774521
$ws.Range("O2").Value = 17.6074171831723

$ws.Range("B3").Value = 15.34727538472584
$ws.Range("C3").Value = 12.22216403410595
$ws.Range("D3").Value = 4.070962645895609
$ws.Range("F3").Value = 19.67037378197129
$ws.Range("G3").Value = 20.80190060518364
$ws.Range("H3").Value = 12.33574090555192
$ws.Range("I3").Value = 19.00274389634443
$ws.Range("L3").Value = 10.59345767434859
$ws.Range("O3").Value = 17.70264518303635

$ws.Range("B4").Value = 14.85359995396853
$ws.Range("C4").Value = 12.10679617159572
$ws.Range("D4").Value = 4.016155791471808
$ws.Range("F4").Value = 19.7103568903177
$ws.Range("G4").Value = 20.86303425462961
$ws.Range("H4").Value = 12.37406142743497
$ws.Range("I4").Value = 19.10136351989577
$ws.Range("L4").Value = 10.57755469158885
$ws.Range("O4").Value = 17.76649962807697

$ws.Range("B5").Value = 14.64740130350778
$ws.Range("C5").Value = 12.05962789835149
$ws.Range("D5").Value = 3.993478285683502
$ws.Range("F5").Value = 19.72850226711406
$ws.Range("G5").Value = 20.89073142647934
$ws.Range("H5").Value = 12.39032500870851
$ws.Range("I5").Value = 19.14284419241717
$ws.Range("L5").Value = 10.57161434807081
$ws.Range("O5").Value = 17.79386939389528

$ws.Range("B6").Value = 14.61286576428441
$ws.Range("C6").Value = 12.05178770451731
$ws.Range("D6").Value = 3.989692548845921
$ws.Range("F6").Value = 19.7316268648973
$ws.Range("G6").Value = 20.8954979703439
$ws.Range("H6").Value = 12.39306466643673
$ws.Range("I6").Value = 19.14981012059033
$ws.Range("L6").Value = 10.57066071954825
$ws.Range("O6").Value = 17.79849541667643

$ws.Range("B7").Value = 14.85083911396882
$ws.Range("C7").Value = 12.106160607028
$ws.Range("D7").Value = 4.015851318592535
$ws.Range("F7").Value = 19.71059411936096
$ws.Range("G7").Value = 20.86339654640063
$ws.Range("H7").Value = 12.37427814210277
$ws.Range("I7").Value = 19.10191770828803
$ws.Range("L7").Value = 10.57747238457478
$ws.Range("O7").Value = 17.76686329308633

$ws.Range("B8").Value = 15.85629808547025
$ws.Range("C8").Value = 12.34460709479525
$ws.Range("D8").Value = 4.128219152552763
$ws.Range("F8").Value = 19.63406399010499
$ws.Range("G8").Value = 20.74625253038363
$ws.Range("H8").Value = 12.29706458475035
$ws.Range("I8").Value = 18.90192164396071
$ws.Range("L8").Value = 10.61226086099594
$ws.Range("O8").Value = 17.63913115623873

$ws.Range("B9").Value = 17.65513741323782
$ws.Range("C9").Value = 12.80373284173752
$ws.Range("D9").Value = 4.33635116113961
$ws.Range("F9").Value = 19.54306570943433
$ws.Range("G9").Value = 20.60614264024513
$ws.Range("H9").Value = 12.16612410371972
$ws.Range("I9").Value = 18.55041431573449
$ws.Range("L9").Value = 10.69713949477577
$ws.Range("O9").Value = 17.43160432607836

$ws.Range("B10").Value = 18.86301971143238
$ws.Range("C10").Value = 13.13321503839212
$ws.Range("D10").Value = 4.480839429995648
$ws.Range("F10").Value = 19.51264527883372
$ws.Range("G10").Value = 20.55901718096955
$ws.Range("H10").Value = 12.0824355761605
$ws.Range("I10").Value = 18.31688574176044
$ws.Range("L10").Value = 10.7693070421408
$ws.Range("O10").Value = 17.30564005508784

$ws.Range("B11").Value = 19.38653059768695
$ws.Range("C11").Value = 13.28086547781843
$ws.Range("D11").Value = 4.54455651513058
$ws.Range("F11").Value = 19.50678498580144
$ws.Range("G11").Value = 20.54990078391902
$ws.Range("H11").Value = 12.04708784652756
$ws.Range("I11").Value = 18.21600122561339
$ws.Range("L11").Value = 10.80418303210672
$ws.Range("O11").Value = 17.25415844920372

$ws.Range("B12").Value = 19.5809543991501
$ws.Range("C12").Value = 13.33641316968988
$ws.Range("D12").Value = 4.568381290470139
$ws.Range("F12").Value = 19.50571668609848
$ws.Range("G12").Value = 20.54823306199669
$ws.Range("H12").Value = 12.03409478585719
$ws.Range("I12").Value = 18.17856727409005
$ws.Range("L12").Value = 10.81767608760842
$ws.Range("O12").Value = 17.23550614243348

$ws.Range("B13").Value = 19.5392527714689
$ws.Range("C13").Value = 13.32446685208895
$ws.Range("D13").Value = 4.563263912877935
$ws.Range("F13").Value = 19.50589553409988
$ws.Range("G13").Value = 20.5485126995455
$ws.Range("H13").Value = 12.03687560997687
$ws.Range("I13").Value = 18.18659515553713
$ws.Range("L13").Value = 10.81475751687489
$ws.Range("O13").Value = 17.23948568261913

$ws.Range("B14").Value = 19.40260292015946
$ws.Range("C14").Value = 13.2854429624154
$ws.Range("D14").Value = 4.546522754263369
$ws.Range("F14").Value = 19.50667401476216
$ws.Range("G14").Value = 20.54972775358404
$ws.Range("H14").Value = 12.0460110312324
$ws.Range("I14").Value = 18.21290610493212
$ws.Range("L14").Value = 10.80528742402182
$ws.Range("O14").Value = 17.25260699203801

$ws.Range("B15").Value = 19.31840140935991
$ws.Range("C15").Value = 13.26149100621389
$ws.Range("D15").Value = 4.536228381252243
$ws.Range("F15").Value = 19.50730081959689
$ws.Range("G15").Value = 20.55070472694126
$ws.Range("H15").Value = 12.05165786265295
$ws.Range("I15").Value = 18.22912243725013
$ws.Range("L15").Value = 10.79952375462192
$ws.Range("O15").Value = 17.26075408065756

$ws.Range("B16").Value = 18.82827698910398
$ws.Range("C16").Value = 13.12351706517641
$ws.Range("D16").Value = 4.476633650243159
$ws.Range("F16").Value = 19.51318914092728
$ws.Range("G16").Value = 20.55986204157333
$ws.Range("H16").Value = 12.08480049135085
$ws.Range("I16").Value = 18.3235863831581
$ws.Range("L16").Value = 10.76706832116799
$ws.Range("O16").Value = 17.30912210724359

$ws.Range("B17").Value = 18.52088648303611
$ws.Range("C17").Value = 13.03827128848008
$ws.Range("D17").Value = 4.439548556953387
$ws.Range("F17").Value = 19.51884779319383
$ws.Range("G17").Value = 20.5686454263464
$ws.Range("H17").Value = 12.10583034258323
$ws.Range("I17").Value = 18.38290664167744
$ws.Range("L17").Value = 10.7476765798525
$ws.Range("O17").Value = 17.34028911082967

$ws.Range("B18").Value = 18.34164402184527
$ws.Range("C18").Value = 12.98903249317449
$ws.Range("D18").Value = 4.418029657562615
$ws.Range("F18").Value = 19.52285329920487
$ws.Range("G18").Value = 20.57485654273323
$ws.Range("H18").Value = 12.11818236867165
$ws.Range("I18").Value = 18.41752955225271
$ws.Range("L18").Value = 10.73671610875261
$ws.Range("O18").Value = 17.35876303149751

$ws.Range("B19").Value = 18.28053946908602
$ws.Range("C19").Value = 12.97232676389615
$ws.Range("D19").Value = 4.410711799677371
$ws.Range("F19").Value = 19.524338314428
$ws.Range("G19").Value = 20.57715817213612
$ws.Range("H19").Value = 12.12240853154669
$ws.Range("I19").Value = 18.42933875768473
$ws.Range("L19").Value = 10.73303849185179
$ws.Range("O19").Value = 17.3651118388144

$ws.Range("B20").Value = 18.55386188191966
$ws.Range("C20").Value = 13.04736766051085
$ws.Range("D20").Value = 4.443515947885595
$ws.Range("F20").Value = 19.51816768492521
$ws.Range("G20").Value = 20.56759035198932
$ws.Range("H20").Value = 12.10356515594841
$ws.Range("I20").Value = 18.37653979468478
$ws.Range("L20").Value = 10.74972093072719
$ws.Range("O20").Value = 17.33691462574196

$ws.Range("B21").Value = 19.44284452629491
$ws.Range("C21").Value = 13.29691544764863
$ws.Range("D21").Value = 4.551448383681305
$ws.Range("F21").Value = 19.50641409926993
$ws.Range("G21").Value = 20.54932234594625
$ws.Range("H21").Value = 12.04331708121248
$ws.Range("I21").Value = 18.2051570776172
$ws.Range("L21").Value = 10.80806131651893
$ws.Range("O21").Value = 17.24873003062708

$ws.Range("B22").Value = 20.00156050183576
$ws.Range("C22").Value = 13.45786710680652
$ws.Range("D22").Value = 4.620213566971679
$ws.Range("F22").Value = 19.50544137144239
$ws.Range("G22").Value = 20.54778821240368
$ws.Range("H22").Value = 12.00622885033451
$ws.Range("I22").Value = 18.09762980124628
$ws.Range("L22").Value = 10.84785474139166
$ws.Range("O22").Value = 17.196010562283

$ws.Range("B23").Value = 19.70542662286219
$ws.Range("C23").Value = 13.37217398132102
$ws.Range("D23").Value = 4.583679134399782
$ws.Range("F23").Value = 19.50534579551309
$ws.Range("G23").Value = 20.54765148856804
$ws.Range("H23").Value = 12.02581395625722
$ws.Range("I23").Value = 18.15460917541485
$ws.Range("L23").Value = 10.82646669056596
$ws.Range("O23").Value = 17.22369640994775

$ws.Range("B24").Value = 18.53896154535559
$ws.Range("C24").Value = 13.04325590337793
$ws.Range("D24").Value = 4.441722906654066
$ws.Range("F24").Value = 19.51847281889602
$ws.Range("G24").Value = 20.5680637343395
$ws.Range("H24").Value = 12.1045884317496
$ws.Range("I24").Value = 18.3794166298085
$ws.Range("L24").Value = 10.74879609318891
$ws.Range("O24").Value = 17.33843850023148

$ws.Range("B25").Value = 17.18796805007744
$ws.Range("C25").Value = 12.68071534334377
$ws.Range("D25").Value = 4.281459299898951
$ws.Range("F25").Value = 19.56131079513816
$ws.Range("G25").Value = 20.63432131975505
$ws.Range("H25").Value = 12.199351828402
$ws.Range("I25").Value = 18.64115931641849
$ws.Range("L25").Value = 10.67242810983474
$ws.Range("O25").Value = 17.48311403341769
